$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @("4-3=1", "11+33=44", "92+1=93", "43-20=23", "8+76=84", "40-39=1", "12+66=78", "98-76=22", "84-26=58", "2+23=25", "87-7=80", "44+0=44", "0+10=10", "71+3=74", "38-7=31", "7+44=51", "78-25=53", "99-65=34", "76-63=13", "83-40=43", "0+1=1", "27+58=85", "73-27=46", "38-24=14", "81-47=34", "65-23=42", "88-44=44", "70-38=32", "67+19=86", "4+89=93", "65-28=37", "80+15=95", "32-28=4", "83-62=21", "56-7=49", "53-48=5", "15-6=9", "40-24=16", "74-42=32", "68+7=75", "82-58=24", "74+17=91", "2+5=7", "72-60=12", "71-64=7", "56+27=83", "5+53=58", "62-38=24", "43+49=92", "68+10=78", "79-76=3", "70+10=80", "27+9=36", "42-25=17", "8+18=26", "31-18=13", "14+82=96", "12+20=32", "46-22=24", "73+17=90", "13+81=94", "21+59=80", "72-17=55", "91-39=52", "60-10=50", "66-19=47", "6+21=27", "79-2=77", "6+54=60", "62-25=37", "99-68=31", "96-51=45", "83-2=81", "5+84=89", "49-18=31", "73-1=72", "67-16=51", "29+38=67", "98-89=9", "28+59=87", "89+10=99", "0+31=31", "34-9=25", "36+8=44", "86-10=76", "58-53=5", "74-51=23", "31+40=71", "25+20=45", "23+5=28", "33+3=36", "93+3=96", "92+0=92", "68+23=91", "67-19=48", "84-64=20", "42-15=27", "2+48=50", "57-26=31", "27+17=44")

$rows = $t.Rows.Count
$cols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $t.Cell($r, $c).Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output "done: $idx cells updated"
